$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "28.501.57"
Set-TextValue $ws.Range("E2") "  +0.94%  "

Set-TextValue $ws.Range("D3") "1.862.90"
Set-TextValue $ws.Range("E3") "  +1.22%  "

Set-TextValue $ws.Range("D4") "1.009"
Set-TextValue $ws.Range("E4") "  +0.37%  "

Set-TextValue $ws.Range("D5") "324.06"
Set-TextValue $ws.Range("E5") "  -0.72%  "

Set-TextValue $ws.Range("D6") "1.008"
Set-TextValue $ws.Range("E6") "  +0.32%  "

Set-TextValue $ws.Range("D7") "0.4543"
Set-TextValue $ws.Range("E7") "  -2.18%  "

Set-TextValue $ws.Range("D8") "0.3821"
Set-TextValue $ws.Range("E8") "  -1.42%  "

Set-TextValue $ws.Range("D9") "0.07784"
Set-TextValue $ws.Range("E9") "  -1.13%  "

Set-TextValue $ws.Range("D10") "0.9813"
Set-TextValue $ws.Range("E10") "  +1.65%  "

Set-TextValue $ws.Range("D11") "21.43"
Set-TextValue $ws.Range("E11") "  -3.21%  "

Set-TextValue $ws.Range("D12") "1.863.31"
Set-TextValue $ws.Range("E12") "  +2.32%  "

Set-TextValue $ws.Range("D13") "6.888"
Set-TextValue $ws.Range("E13") "  -0.17%  "

Set-TextValue $ws.Range("D14") "5.610"
Set-TextValue $ws.Range("E14") "  -1.69%  "

Set-TextValue $ws.Range("D15") "0.06950"
Set-TextValue $ws.Range("E15") "  +1.10%  "

Set-TextValue $ws.Range("D16") "87.24"
Set-TextValue $ws.Range("E16") "  -1.63%  "

Set-TextValue $ws.Range("D17") "1.009"
Set-TextValue $ws.Range("E17") "  +0.43%  "

Set-TextValue $ws.Range("D18") "0.000009933"
Set-TextValue $ws.Range("E18") "  -0.43%  "

Set-TextValue $ws.Range("D19") "16.59"
Set-TextValue $ws.Range("E19") "  -1.01%  "

Set-TextValue $ws.Range("D20") "1.007"
Set-TextValue $ws.Range("E20") "  +0.38%  "

Set-TextValue $ws.Range("D21") "28.519.92"
Set-TextValue $ws.Range("E21") "  +0.98%  "

Set-TextValue $ws.Range("D22") "5.229"
Set-TextValue $ws.Range("E22") "  -1.60%  "

Set-TextValue $ws.Range("D23") "10.86"
Set-TextValue $ws.Range("E23") "  -1.84%  "

Set-TextValue $ws.Range("D24") "2.109"
Set-TextValue $ws.Range("E24") "  -0.15%  "

Set-TextValue $ws.Range("D25") "2.129.92"
Set-TextValue $ws.Range("E25") "  +3.77%  "

Set-TextValue $ws.Range("D26") "152.55"
Set-TextValue $ws.Range("E26") "  -1.38%  "

Set-TextValue $ws.Range("D27") "19.09"
Set-TextValue $ws.Range("E27") "  -0.60%  "

Set-TextValue $ws.Range("D28") "5.614"
Set-TextValue $ws.Range("E28") "  -2.32%  "

Set-TextValue $ws.Range("D29") "1.929"
Set-TextValue $ws.Range("E29") "  -2.42%  "

Set-TextValue $ws.Range("D30") "117.69"
Set-TextValue $ws.Range("E30") "  -1.31%  "

Set-TextValue $ws.Range("D31") "0.09252"
Set-TextValue $ws.Range("E31") "  -0.19%  "

Set-TextValue $ws.Range("D32") "0.9007"
Set-TextValue $ws.Range("E32") "  -4.00%  "

Set-TextValue $ws.Range("D33") "5.262"
Set-TextValue $ws.Range("E33") "  -0.61%  "

Set-TextValue $ws.Range("E34") "  -0.95%  "

Set-TextValue $ws.Range("D35") "3.308"
Set-TextValue $ws.Range("E35") "  -0.72%  "

Set-TextValue $ws.Range("D36") "0.05679"
Set-TextValue $ws.Range("E36") "  -2.68%  "

Set-TextValue $ws.Range("D37") "1.134"
Set-TextValue $ws.Range("E37") "  -0.70%  "

Set-TextValue $ws.Range("D38") "0.02055"
Set-TextValue $ws.Range("E38") "  -3.41%  "

Set-TextValue $ws.Range("D39") "7.681"
Set-TextValue $ws.Range("E39") "  -1.19%  "

Set-TextValue $ws.Range("D40") "0.5542"
Set-TextValue $ws.Range("E40") "  -1.12%  "

Set-TextValue $ws.Range("D41") "0.1768"
Set-TextValue $ws.Range("E41") "  +0.18%  "

Set-TextValue $ws.Range("D42") "9.598"
Set-TextValue $ws.Range("E42") "  -3.44%  "

Set-TextValue $ws.Range("D43") "0.07087"
Set-TextValue $ws.Range("E43") "  -3.17%  "

Set-TextValue $ws.Range("D44") "11.60"
Set-TextValue $ws.Range("E44") "  -0.21%  "

Set-TextValue $ws.Range("D45") "0.5216"
Set-TextValue $ws.Range("E45") "  -1.26%  "

Set-TextValue $ws.Range("D46") "2.121"
Set-TextValue $ws.Range("E46") "  -0.92%  "

Set-TextValue $ws.Range("B47") "WEMIXToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D47") "1.106"
Set-TextValue $ws.Range("E47") "  -3.16%  "

Set-TextValue $ws.Range("B48") "NEARProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.805"
Set-TextValue $ws.Range("E48") "  -1.89%  "

Set-TextValue $ws.Range("B49") "Quant"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D49") "112.03"
Set-TextValue $ws.Range("E49") "  -1.83%  "

Set-TextValue $ws.Range("D50") "2.420"
Set-TextValue $ws.Range("E50") "  +4.09%  "

Set-TextValue $ws.Range("D51") "1.007"
Set-TextValue $ws.Range("E51") "  +0.27%  "
